$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Logs" sheet: append row 4 with the new test-mail entry
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A4").Value = "Hoe kan ik een product retourneren?"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Testmail #3: Hoe kan ik een product retourneren?"
$logs.Range("D4").Value = "Retour / Terugbetaling"
$logs.Range("E4").Value = "Geachte klant,`nBedankt voor uw bericht. Als u een product wilt retourneren, kunt u dit doen door contact op te nemen met onze klantenservice. Zij zullen u voorzien van verdere instructies met betrekking tot de retourprocedure en eventuele benodigde informatie.`nVoor een snelle afhandeling van uw verzoek, verzoeken wij u vriendelijk om uw ordernummer en reden van retournering te vermelden in uw bericht naar onze klantenservice.`nMet vriendelijke groet,`n[Naam Bedrijf] E-mailassistent"
$logs.Range("F4").Value = "2025-07-22 12:16:09"
$logs.Range("G4").Value = "Ja"
$logs.Range("H4").Value = "Nee"
$logs.Range("I4").Value = "Ja"
$logs.Range("J4").Value = "Nee"

# The multi-line Antwoord text (embedded newlines) makes the host auto-size
# the row; reset it back to the sheet's normal (non-custom) row height so the
# new row matches the existing rows 2-3, which carry no explicit height.
$logs.Rows.Item(4).AutoFit()

# Extend the conditional-formatting ranges so they cover the new row too.
$ranges = @("D2:D3", "G2:G3", "H2:H3", "I2:I3", "J2:J3")
foreach ($rngAddr in $ranges) {
    $col = $rngAddr.Substring(0, 1)
    $newAddr = "$col" + "2:" + "$col" + "4"
    $fcs = $logs.Range($rngAddr).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newAddr))
    }
}

# ---------------------------------------------------------------------------
# 2. "Dashboard" sheet: append row 4 with the new category tally
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Retour / Terugbetaling"
$dash.Range("B4").Value = 1

# ---------------------------------------------------------------------------
# 3. Chart on the Dashboard sheet: extend the category / value series refs
# ---------------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$4,Dashboard!`$B`$2:`$B`$4,1)"
